$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kinetics1")
$c = $ws.Columns.Item(2)
Write-Host ("ColumnWidth: " + $c.ColumnWidth)
$c.ColumnWidth = 31.61
Write-Host ("ColumnWidth after: " + $c.ColumnWidth)
